# Add a new "Crow" translation entry as a new row inserted at row 50,
# pushing the existing rows 50-130 down to 51-131.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50").Insert()

$ws.Range("A50").Value = "object.CROW"
$ws.Range("B50").Value = "Crow"
$ws.Range("C50").Value = "Cuervo"

# Update the view state to match the post-edit selection/zoom.
$win = $excel.ActiveWindow
$win.Zoom = 115
[void]$ws.Range("D50").Select()
